# Daily attendance processing - 2025-10-30 23:20:54
#
# Column G ("Recorded By") holds a comma-separated list of recorder
# identities, e.g. "System, dnasr281@gmail.com". For every row whose list
# starts with the literal token "System" (and has more than one entry),
# flip the order of the list so "System" moves from first to last.
# Rows that don't start with "System", or only contain a single entry,
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }
    if ($val.Length -eq 0) {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Length -gt 1 -and $parts[0] -ceq "System") {
        $n = $parts.Length
        $rev = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $rev += $parts[$i]
        }
        $newVal = [string]::Join(", ", $rev)
        $cell.Value = $newVal
    }
}
